# edit.ps1 - applies the recorded change to before.pptx:
#   1. Three tables (on slides 14, 15, 16) switch from the custom
#      "Table_0" style {F63DF6F2-5E40-4309-9412-EA3D2AF853E3} to the
#      built-in table style {50C7EF46-2E12-470A-9EEE-DF136F7E0619}.
#   2. The deck's theme (ppt/theme/theme1.xml, used by the slide master -
#      and therefore by every slide) is swapped from the "Integral" /
#      "Red Violet" palette to the default "Office Theme" / "Office"
#      palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newTableStyle = "{50C7EF46-2E12-470A-9EEE-DF136F7E0619}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Swap the theme colors used by the slide master/theme1.xml ----
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (Theme color index 1-12)
$officePalette = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $rrggbb = $officePalette[$i - 1]
    $r = [math]::Floor($rrggbb / 0x10000) -band 0xFF
    $g = [math]::Floor($rrggbb / 0x100) -band 0xFF
    $b = $rrggbb -band 0xFF
    $bgr = ($b * 0x10000) + ($g * 0x100) + $r
    $themeColors.Colors($i).RGB = $bgr
}
